$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 61.91587448120117
$ws.Range("E3").Value = 113.3610486984253
$ws.Range("E4").Value = 159.392237663269
$ws.Range("E5").Value = 195.9049701690674
$ws.Range("E6").Value = 208.9052796363831
$ws.Range("E7").Value = 228.295361995697
$ws.Range("E8").Value = 865167.9310798645
$ws.Range("E11").Value = 47.99246788024902
$ws.Range("E12").Value = 67.71337985992432
$ws.Range("E13").Value = 83.60910415649414
$ws.Range("E14").Value = 93.8071608543396
$ws.Range("E15").Value = 293274.8041152954
$ws.Range("E16").Value = 10.19382476806641
$ws.Range("E17").Value = 14.39611911773682
$ws.Range("E18").Value = 18.35793256759644
$ws.Range("E19").Value = 22.87697792053223
$ws.Range("E20").Value = 23.88185262680054
$ws.Range("E21").Value = 25.51437616348267
$ws.Range("E22").Value = 293192.8429603577
$ws.Range("E23").Value = 9.02104377746582
$ws.Range("E24").Value = 13.05216550827026
$ws.Range("C25").Value = 0.8338607594936709
$ws.Range("E25").Value = 16.32720232009888
$ws.Range("F25").Value = 0.9738529014844804
$ws.Range("G25").Value = 0.7848075348075347
$ws.Range("C26").Value = 0.9705882352941176
$ws.Range("E26").Value = 20.3709602355957
$ws.Range("G26").Value = 0.9629629629629629
$ws.Range("C27").Value = 0.9934640522875817
$ws.Range("E27").Value = 21.12871408462524
$ws.Range("E28").Value = 22.62958288192749
$ws.Range("E29").Value = 995372.3289966583
$ws.Range("E33").Value = 1.961946487426758
$ws.Range("E34").Value = 2.634048461914062
$ws.Range("E35").Value = 4.368138313293456
$ws.Range("E36").Value = 90.25192260742188
$ws.Range("E38").Value = 60.20407676696777
$ws.Range("E39").Value = 85.50363779067993
$ws.Range("E40").Value = 119.3625926971436
$ws.Range("E41").Value = 136.11900806427
$ws.Range("E42").Value = 145.0353503227234
$ws.Range("E43").Value = 212.2581005096436
$ws.Range("E44").Value = 4.981040954589844
$ws.Range("E45").Value = 7.868587970733643
$ws.Range("E46").Value = 10.69271564483643
$ws.Range("E47").Value = 16.26861095428467
$ws.Range("E48").Value = 18.64022016525269
$ws.Range("E49").Value = 20.59328556060791
$ws.Range("E50").Value = 40.16709327697754
$ws.Range("E52").Value = 1.489649998955429
$ws.Range("E53").Value = 2.518249998800457
$ws.Range("E54").Value = 2.885500002652407
$ws.Range("E55").Value = 3.099999996833503
$ws.Range("E56").Value = 3.544799996167421
$ws.Range("E57").Value = 94.54900000244379
$ws.Range("E58").Value = 1.650094985961914
$ws.Range("E59").Value = 1.862752437591553
$ws.Range("E60").Value = 1.988232135772705
$ws.Range("E61").Value = 2.152919769287109
$ws.Range("E62").Value = 2.233207225799561
$ws.Range("E63").Value = 2.337312698364258
$ws.Range("E64").Value = 60.42885780334473
$ws.Range("E65").Value = 1.717090606689453
$ws.Range("E66").Value = 2.348566055297852
$ws.Range("E67").Value = 2.612411975860596
$ws.Range("E68").Value = 3.137588500976562
$ws.Range("E69").Value = 3.612220287322998
$ws.Range("E70").Value = 4.477226734161376
$ws.Range("E71").Value = 22.4609375
$ws.Range("E72").Value = 6.799221038818359
$ws.Range("E73").Value = 12.29087114334106
$ws.Range("E74").Value = 17.89242029190063
$ws.Range("E75").Value = 25.15101432800293
$ws.Range("E76").Value = 25.39414167404175
$ws.Range("E77").Value = 25.90879201889038
$ws.Range("E78").Value = 30.94387054443359
$ws.Range("E79").Value = 2.050161361694336
$ws.Range("E80").Value = 2.219116687774658
$ws.Range("E81").Value = 2.424299716949463
$ws.Range("E82").Value = 2.60007381439209
$ws.Range("E83").Value = 2.744913101196289
$ws.Range("E84").Value = 3.067731857299804
$ws.Range("E85").Value = 8.311033248901367
$ws.Range("E86").Value = 1.193046569824219
$ws.Range("E87").Value = 1.302719116210938
$ws.Range("E88").Value = 1.600265502929688
$ws.Range("E89").Value = 2.04002857208252
$ws.Range("E90").Value = 2.304792404174805
$ws.Range("E91").Value = 2.665376663208008
$ws.Range("E92").Value = 9.670734405517578
$ws.Range("E94").Value = 1.716300001367927
$ws.Range("E95").Value = 2.034250000491738
$ws.Range("E96").Value = 2.268499996513128
$ws.Range("E97").Value = 2.668499996885657
$ws.Range("E98").Value = 4.646150003001093
$ws.Range("E99").Value = 32.17499999701977
$ws.Range("E101").Value = 0.6296499971300363
$ws.Range("E102").Value = 0.8002500031143427
$ws.Range("E103").Value = 0.96000000461936
$ws.Range("E104").Value = 1.33124999050051
$ws.Range("E105").Value = 1.815750003606081
$ws.Range("E106").Value = 9.525999993085861
$ws.Range("E107").Value = 2.655982971191406
$ws.Range("E108").Value = 4.889726638793945
$ws.Range("E109").Value = 6.480157375335693
$ws.Range("E110").Value = 8.646011352539062
$ws.Range("E111").Value = 8.827269077301025
$ws.Range("E112").Value = 8.926045894622803
$ws.Range("E113").Value = 13.62109184265137
$ws.Range("E114").Value = 8.629000000655651
$ws.Range("E115").Value = 12.3289999961853
$ws.Range("E116").Value = 17.54925000108778
$ws.Range("E117").Value = 36.60999999940395
$ws.Range("E118").Value = 39.41549999639392
$ws.Range("E119").Value = 45.6207500046119
$ws.Range("E120").Value = 154.875
$ws.Range("E121").Value = 63.06099891662598
$ws.Range("E122").Value = 126.7129898071289
$ws.Range("E123").Value = 182.6382875442505
$ws.Range("E124").Value = 253.0834674835205
$ws.Range("E125").Value = 254.905104637146
$ws.Range("E126").Value = 259.7604274749756
$ws.Range("E127").Value = 293460.0381851196
